$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New localization rows to append at rows 60-62 (id / filter / text_JP / text columns: A / B / C / D)
# Row 60
$ws.Range("A60").Value = "cwl_warn_missing_mods"
$ws.Range("C60").Value = "現在のセーブから欠落しているMOD：`n{0}"
$ws.Range("D60").Value = "当前存档中缺失的模组：`n{0}"

# Row 61
$ws.Range("A61").Value = "cwl_warn_missing_mods_yes"
$ws.Range("C61").Value = "セーブせずに終了"
$ws.Range("D61").Value = "不保存并返回至标题"

# Row 62
$ws.Range("A62").Value = "cwl_warn_missing_mods_no"
$ws.Range("C62").Value = "プレイを続ける"
$ws.Range("D62").Value = "继续游玩"

# Copy style from existing formatted rows so formatting matches the target layout
$ws.Range("A60").Style = $ws.Range("A59").Style
$ws.Range("C60").Style = $ws.Range("C59").Style
$ws.Range("D60").Style = $ws.Range("D58").Style

$ws.Range("A61").Style = $ws.Range("A59").Style
$ws.Range("C61").Style = $ws.Range("C59").Style
$ws.Range("D61").Style = $ws.Range("D59").Style

$ws.Range("A62").Style = $ws.Range("A59").Style
$ws.Range("C62").Style = $ws.Range("C59").Style
$ws.Range("D62").Style = $ws.Range("D59").Style

# D59 style changes from s="10" to s="3" in the diff
$ws.Range("D59").Style = $ws.Range("C59").Style

# Row heights: row 60 becomes 46.5pt (matching rows with wrapped multi-line rich text), rows 61/62 stay 23.25pt
$ws.Rows.Item(60).RowHeight = 46.5
$ws.Rows.Item(61).RowHeight = 23.25
$ws.Rows.Item(62).RowHeight = 23.25

# Selection / view state to mirror diff (sheetView topLeftCell/selection)
$ws.Range("D65").Select()
$excel.ActiveWindow.ScrollRow = 53
